# Trade #25 closed at 2026-02-17 08:03:13 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the
# "MarketMaking" strategy and appends the newly-closed trade (#25) to both
# the "All Trades" and "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.54   # Current Capital
$summary.Range("B4").Value = -0.46     # Total P&L $
$summary.Range("B5").Value = -0.37     # Total P&L %
$summary.Range("B6").Value = 25        # Total Trades
$summary.Range("B8").Value = 12        # Losing Trades
$summary.Range("B9").Value = 28        # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.54      # Capital
$status.Range("D4").Value = 25         # Trades
$status.Range("E4").Value = -0.46      # P&L $
$status.Range("F4").Value = -0.46      # P&L %
$status.Range("G4").Value = 28         # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new trade row (#25, spreadsheet row 26) to both the
#    "All Trades" and "MarketMaking" logs - they share identical content.
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A26").Value = 25

    # Date/time columns are stored as plain text in this workbook (not
    # real Excel dates), so force text with a leading apostrophe and then
    # strip the style edit it would otherwise pick up.
    $ws.Range("B26").Value = "'2026-02-17"
    $ws.Range("B26").Style = "Normal"
    $ws.Range("C26").Value = "08:03:07"

    $ws.Range("D26").Value = "MarketMaking"
    $ws.Range("E26").Value = "UP"
    $ws.Range("F26").Value = 0.34
    $ws.Range("G26").Value = 0.26
    $ws.Range("H26").Value = "CLOSED"
    $ws.Range("I26").Value = -23.5294
    $ws.Range("J26").Value = -0.08
    $ws.Range("K26").Value = 99.54
    $ws.Range("L26").Value = 0
    $ws.Range("M26").Value = 0
    $ws.Range("N26").Value = 0.6
    $ws.Range("O26").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P26").Value = "early_exit"
    $ws.Range("Q26").Value = 0.13
}
